$wb = $excel.ActiveWorkbook

# --- Update arrivingData sheet content first (append city1..city5 in column B) ---
# (done before Searchflightdata edits so new shared strings land in the same
# order as the target workbook: city1..city5 then fromcity/tocity)
$wsArrive = $wb.Worksheets.Item("arrivingData")
$wsArrive.Range("B1").Value = "city1"
$wsArrive.Range("B2").Value = "city2"
$wsArrive.Range("B3").Value = "city3"
$wsArrive.Range("B4").Value = "city4"
$wsArrive.Range("B5").Value = "city5"

# --- Rename sheet "departingData" -> "Searchflightdata" ---
$wsSearch = $wb.Worksheets.Item("departingData")
$wsSearch.Name = "Searchflightdata"

# --- Update Searchflightdata (former departingData) sheet content ---
$wsSearch.Range("A1").Value = "fromcity"
$wsSearch.Range("A2").Value = "tocity"
$wsSearch.Range("A3").ClearContents()
$wsSearch.Range("A4").ClearContents()
$wsSearch.Range("A5").ClearContents()
$wsSearch.Range("B1").Value = "Frankfurt"
$wsSearch.Range("B2").Value = "London"

# Match formatting (style index 1 - Consolas 9pt FF222222) already used by A2:A5
$wsSearch.Range("A2").Copy()
$wsSearch.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$wsSearch.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$wsSearch.Range("B2").PasteSpecial(-4122)  # xlPasteFormats

# --- Selections on each sheet ---
$wsArrive.Activate()
$wsArrive.Range("C3").Select()

$wsSearch.Activate()
$wsSearch.Range("G22").Select()

# --- Window sizing (best effort) ---
$win = $excel.ActiveWindow
$win.Height = 3765
